# Apply updated power-flow results ("case with 380 kV done") to Sheet1.
# Columns D, J, K are untouched (they remain 0 in every data row);
# columns B, C, E, F, G, H, I, L, M, N, O receive new values for rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cols = @('B','C','E','F','G','H','I','L','M','N','O')

$data = @{
    2  = @(0.6742382825662219, 0.281511994406145, 0.134278848984164, 0.4443680307746121, 0.3257080813141897, 0.4987926123142614, 0.4851656965642519, 0.1987929869517586, 0.1659874430461628, 1.167285015218582, 1.581121824911321)
    3  = @(0.6017064188764891, 0.2741861385872397, 0.1355714005859827, 0.3878228170618172, 0.3263827970977715, 0.5024417969223833, 0.4912712563549189, 0.1960395792534655, 0.1537947441846441, 1.169665638044904, 1.589961888529501)
    4  = @(0.5571227097429698, 0.2696687952539492, 0.1364172345828063, 0.3531389305169483, 0.3270662036578855, 0.504920064490058, 0.4953065582682274, 0.1944455823451676, 0.1463380904184781, 1.171540430084612, 1.596449701326023)
    5  = @(0.5389436534278786, 0.267823271135029, 0.1367750596031918, 0.3390132514313251, 0.3274122923205098, 0.5059897726737788, 0.4970229793268199, 0.193820389404479, 0.1433071606388907, 1.172408570218188, 1.599359975097713)
    6  = @(0.5359244234243477, 0.2675165461388502, 0.1368352702631022, 0.336668177824194, 0.3274738403063395, 0.5061710093196652, 0.4973123376200235, 0.1937180510951464, 0.142804350387248, 1.172559022601824, 1.599859312981863)
    7  = @(0.5568775825702232, 0.2696439245556377, 0.1364220071104905, 0.3529483938344953, 0.3270705975470136, 0.5049342487995574, 0.495329415027042, 0.1944370519780918, 0.1462971826672437, 1.17155171607083, 1.596487871637635)
    8  = @(0.6492402830622837, 0.2789901304630007, 0.1347136948324898, 0.4248636149813478, 0.3258848170991655, 0.5000015465141558, 0.4872114200023443, 0.1978236087177834, 0.1617773770629398, 1.168020281436497, 1.583949816566061)
    9  = @(0.8299204910541107, 0.2971588410880202, 0.1317772417785079, 0.5661985755041457, 0.3256989697415378, 0.4922128275124749, 0.4735663138110411, 0.2052281445628452, 0.1923604831164951, 1.164361513130288, 1.567779131830434)
    10 = @(0.9623348603564637, 0.3104028226835283, 0.1298709984861528, 0.6702781546542269, 0.3268729759421305, 0.4876374993081143, 0.4649293672356265, 0.2111304593802714, 0.2149577669725886, 1.163650791764681, 1.561039539141149)
    11 = @(1.022489363041188, 0.3164036256189604, 0.1290581376239919, 0.7176906081379002, 0.3276930106115685, 0.4858047970777619, 0.461301967304717, 0.2139152952961467, 0.2252635946075685, 1.163754032107931, 1.559092131026773)
    12 = @(1.045255278938782, 0.3186723746819666, 0.1287581220293621, 0.7356546913071611, 0.3280447554990786, 0.4851465260974521, 0.4599717694298882, 0.2149841259661684, 0.2291696910547785, 1.16385422067799, 1.558515683628855)
    13 = @(1.040352847691622, 0.3181839230617811, 0.1288223891481408, 0.7317853510981394, 0.3279671664062107, 0.4852867075932465, 0.4602563194230029, 0.2147533006625935, 0.2283282917205511, 1.163829929962631, 1.558632669139371)
    14 = @(1.02436260398423, 0.316590350647374, 0.1290332989559975, 0.7191683204515869, 0.327721122326551, 0.4857499246893582, 0.4611916606353184, 0.2140029430510708, 0.2255848828257783, 1.16376105165719, 1.559041478893448)
    15 = @(1.014566340707688, 0.3156137641239241, 0.1291635024015436, 0.7114413442032514, 0.3275757834399258, 0.4860383115024831, 0.4617702401883186, 0.2135451837256852, 0.223904913996968, 1.163726810733749, 1.559312857501624)
    16 = @(0.9584018482861438, 0.3100101586907442, 0.1299252125508072, 0.6671810134426437, 0.3268251472388144, 0.4877622722076396, 0.4651724995266164, 0.2109504657655492, 0.2142847615875922, 1.163652607136711, 1.561189326966741)
    17 = @(0.9239247664736467, 0.3065662651808623, 0.1304063964921331, 0.6400460337125793, 0.3264379674397588, 0.4888835339522117, 0.4673369566792047, 0.2093842032004147, 0.2083896335142867, 1.163716160943267, 1.562627060889113)
    18 = @(0.9040868876838317, 0.3045831851354706, 0.1306882722781042, 0.6244449056556647, 0.326242183556829, 0.4895518582035692, 0.4686102792513225, 0.2084927331367794, 0.2050013946394031, 1.163792878758727, 1.563559273416189)
    19 = @(0.8973688701170772, 0.3039113684804136, 0.1307845888158621, 0.619163680173358, 0.3261805137544016, 0.4897821616114655, 0.4690462778130495, 0.2081925141630165, 0.2038546297106905, 1.163825759676158, 1.563892979182114)
    20 = @(0.9275957059049347, 0.306933107038418, 0.130354644697519, 0.642933953830422, 0.3264763974532414, 0.4887617516292408, 0.4671036085023452, 0.2095499619815087, 0.2090169249190765, 1.163705240554279, 1.562463115856218)
    21 = @(1.029059699364609, 0.3170585213073025, 0.1289711381068201, 0.7228739723491628, 0.327792272124924, 0.4856128969622659, 0.4609157494408329, 0.2142229545524543, 0.2263905955611349, 1.163779626745239, 1.558917030976858)
    22 = @(1.095294098208285, 0.3236548534239034, 0.1281123794076171, 0.7751780083420101, 0.3288925502356363, 0.4837632116275898, 0.4571247329673298, 0.217360184687422, 0.2377655896976734, 1.164184232611277, 1.557537926399192)
    23 = @(1.059951226920305, 0.3201362672477899, 0.1285665606677231, 0.7472568307830727, 0.3282832953041463, 0.4847313731822283, 0.4591248936005314, 0.2156782045205148, 0.231692771931904, 1.163935793348045, 1.558188058739717)
    24 = @(0.9259361256035277, 0.3067672676015718, 0.1303780253714919, 0.6416283278902171, 0.3264589397378685, 0.4888167356050133, 0.4672090150487698, 0.2094749944119343, 0.2087333234994304, 1.163710052486664, 1.562536906358417)
    25 = @(0.7810954713378351, 0.2922615446413062, 0.1325274546587494, 0.5279251897347166, 0.3255195491309522, 0.4941183286039958, 0.4770140286164022, 0.2031436181353996, 0.223904913996968, 1.165003070533714, 1.571251638662744)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
